$wb = $excel.ActiveWorkbook

# The existing "InvalidNegotiationsFee" sheet is the structural template for
# these "broken recalculation" test fixtures (same headers/layout/styles).
# Duplicate it and drop the copy right after the original, then rename it
# for the new scenario being covered: a 'ServiceTax' ('ISS') that doesn't
# recalculate as 'Brokerage' ('Corretagem') * 'ServiceTaxRate'.
$template = $wb.Worksheets.Item("InvalidNegotiationsFee")
$template.Copy($null, $template) | Out-Null
$newSheet = $wb.Worksheets.Item($template.Index + 1)
$newSheet.Name = "InvalidServiceTax"

# Make the settlement fee (H) a proper formula again (it was a hard-coded
# literal on the template sheet).
$newSheet.Range("H2").Formula = "=F2*0.005%"

# Break the service tax (J) for the first line: instead of the expected
# 'Brokerage' * 'ServiceTaxRate' formula, it's now a hard-coded value that
# doesn't match what recalculation would produce.
$newSheet.Range("J2").Value = 0.12

# The second line's service tax keeps being a correctly-calculated formula.
$newSheet.Range("J3").Formula = "=I3*6.5%"

# Match the author's on-disk selection/active-cell state for the new sheet.
$newSheet.Activate() | Out-Null
$newSheet.Range("J2").Select() | Out-Null

$wb.Application.CalculateFull() | Out-Null
